$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.246.70"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.658.37"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "218.31"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "0.5332"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").Value = "0.2631"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "0.06354"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "20.51"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("D11").Value = "0.07840"
$ws.Range("D12").Value = "4.533"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").Value = "1.684.60"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "1.886.53"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "0.5508"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "0.0₅8180"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "65.52"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "26.228.56"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "4.631"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").Value = "191.87"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").Value = "10.12"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").Value = "6.038"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").Value = "144.31"
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("D26").Value = "0.1227"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").Value = "7.224"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("E29").Value = "  +2.71%  "
$ws.Range("D30").Value = "0.05789"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").Value = "1.278"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "3.569"
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("D33").Value = "3.280"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("D36").Value = "0.9556"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").Value = "0.5785"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").Value = "0.01603"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "5.834"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "0.8523"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").Value = "104.74"
$ws.Range("E43").Value = "  +4.13%  "
$ws.Range("D44").Value = "1.042.04"
$ws.Range("E44").Value = "  +3.39%  "
$ws.Range("D45").Value = "1.799.60"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "56.95"
$ws.Range("D47").Value = "1.010"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.4367"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₈103"
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("D50").Value = "7.931"
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("D51").Value = "0.05156"
$ws.Range("E51").Value = "  +0.08%  "
